# Femacal de La Calera - Poroto granado: weekly data refresh.
# Two new daily observations are inserted into the historical series:
#   - one at (what becomes) row 23, pushing the former rows 23..104 down by one
#   - a second one further down at (what becomes) row 67, pushing the rows
#     that are by then at 67..105 down by one more
# so the sheet grows from A1:R104 to A1:R106.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DataRow($Row, $Fecha, $Volumen, $PrecioMinimo, $PrecioMaximo, $PrecioPromedio, $Unidad, $Origen, $PrecioKg) {
    $ws.Cells.Item($Row, 1).Value = 3
    $ws.Cells.Item($Row, 2).Value = "Femacal de La Calera"
    $ws.Cells.Item($Row, 3).Value = "Coquimbo"
    $ws.Cells.Item($Row, 4).Value = $Fecha
    $ws.Cells.Item($Row, 5).Value = 5
    $ws.Cells.Item($Row, 6).Value = 100112030
    $ws.Cells.Item($Row, 7).Value = "Poroto granado"
    $ws.Cells.Item($Row, 8).Value = "Sin especificar"
    $ws.Cells.Item($Row, 9).Value = "Primera"
    $ws.Cells.Item($Row, 10).Value = $Volumen
    $ws.Cells.Item($Row, 11).Value = $PrecioMinimo
    $ws.Cells.Item($Row, 12).Value = $PrecioMaximo
    $ws.Cells.Item($Row, 13).Value = $PrecioPromedio
    $ws.Cells.Item($Row, 14).Value = $Unidad
    $ws.Cells.Item($Row, 15).Value = $Origen
    $ws.Cells.Item($Row, 16).Value = $PrecioKg
    $ws.Cells.Item($Row, 17).Value = 25
    $ws.Cells.Item($Row, 18).Value = "Hortaliza"
}

# Insert the first new record before the old row 23 -- everything from the
# old row 23 through the old last row (104) shifts down by one.
$ws.Rows(23).Insert()
Set-DataRow 23 44188 45 40000 40000 40000 "$/saco 25 kilos" "Provincia de Talca" 1600

# Insert the second new record before (what is now) row 67 -- everything
# from there through the new last row (105) shifts down by one more.
$ws.Rows(67).Insert()
Set-DataRow 67 44518 65 37000 38000 37538 "$/malla 25 kilos" "Provincia de Limarí" 1502
